$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the two "section header only" rows (originally row 5 "situação do domicílio"
# and row 8 "grandes regiões e unidades da federação") which had no numeric data.
# Delete from the bottom up so row numbers of earlier rows stay valid.
$ws.Rows(8).Delete() | Out-Null
$ws.Rows(5).Delete() | Out-Null

# Fix sub-header labels on row 2 (B2 and F2 should read "total" instead of the
# stray "unnamed: 1_level_1" / "unnamed: 5_level_1" placeholders)
$ws.Range("B2").Value = "total"
$ws.Range("F2").Value = "total"
